# Fixed naive component forecaster bug - Presentation state 11.02.
# Row 2 and Row 3: the C column (y_0_forecast) value for the very first
# two rows was bogus data leftover from a prior run; clear it back out
# (and the stray E2 value that went with it). Remaining C/E cells get
# refreshed with the corrected (recomputed) forecast values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()

$ws.Range("E3").Value = 0.6970543652217387
$ws.Range("C4").Value = -0.01587181126745385
$ws.Range("C5").Value = -0.02256889165886955
$ws.Range("E5").Value = -0.02753509623224515
$ws.Range("C6").Value = 0.09611428386595566
$ws.Range("E7").Value = -0.2251688766574889
$ws.Range("C8").Value = -0.001350220946472191
$ws.Range("C10").Value = -0.5761528471665334
$ws.Range("C14").Value = -0.4278219446121501
$ws.Range("C15").Value = -1.026566979837429
$ws.Range("C17").Value = 0.4636049209196802
$ws.Range("E17").Value = 0.2986939435938973
$ws.Range("C18").Value = 0.6216390921348403
$ws.Range("E18").Value = -0.0776179936130994
$ws.Range("C19").Value = -0.6768900623516871
